$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "117-65549061"
$ws.Range("B2").Value = 2870326740
$ws.Range("C2").Value = 2870326740

# Row 3
$ws.Range("A3").Value = "117-36803410"
$ws.Range("B3").Value = 2870325422
$ws.Range("C3").Value = 2870325422

# Row 4
$ws.Range("A4").Value = "117-33837241"
$ws.Range("B4").Value = 2880048980
$ws.Range("C4").Value = 2880048980

# Row 5
$ws.Range("A5").Value = "117-37204893"
$ws.Range("B5").Value = 2870326149
$ws.Range("C5").Value = 2870326149

# Row 6
$ws.Range("A6").Value = "117-89697893"
$ws.Range("B6").Value = 2870326380
$ws.Range("C6").Value = 2870326380

# Row 7
$ws.Range("A7").Value = "117-89666640"
$ws.Range("B7").Value = 216005522
$ws.Range("C7").Value = 216005522

# Row 8
$ws.Range("A8").Value = "117-89789055"
$ws.Range("B8").Value = 216005530
$ws.Range("C8").Value = 216005530

# Row 9
$ws.Range("A9").Value = "117-89801574"
$ws.Range("B9").Value = 2880048998
$ws.Range("C9").Value = 2880048998

# Row 10
$ws.Range("A10").Value = "117-35547223"
$ws.Range("B10").Value = 220427831
$ws.Range("C10").Value = 220427831

# Row 11
$ws.Range("A11").Value = "117-36885774"
$ws.Range("B11").Value = 2870326840
$ws.Range("C11").Value = 2870326840
